$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 146, shifting rows 146:264 down to 147:265
$ws.Rows(146).Insert()

# Populate the newly inserted row 146 with the new data record
$ws.Cells.Item(146, 1).Value = 6
$ws.Cells.Item(146, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(146, 3).Value = "Metropolitana"

[datetime]$fecha146 = "2023-01-16"
$ws.Cells.Item(146, 4).Value = $fecha146

$ws.Cells.Item(146, 5).Value = 13
$ws.Cells.Item(146, 6).Value = "Fruta"
$ws.Cells.Item(146, 7).Value = 100101
$ws.Cells.Item(146, 8).Value = "Berries"
$ws.Cells.Item(146, 9).Value = 100101004
$ws.Cells.Item(146, 10).Value = "Frambuesa"
$ws.Cells.Item(146, 11).Value = "Sin especificar"
$ws.Cells.Item(146, 12).Value = "Especial"
$ws.Cells.Item(146, 13).Value = 70
$ws.Cells.Item(146, 14).Value = 7000
$ws.Cells.Item(146, 15).Value = 7000
$ws.Cells.Item(146, 16).Value = 7000
$ws.Cells.Item(146, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(146, 18).Value = "Región del Maule"
$ws.Cells.Item(146, 19).Value = 3500
$ws.Cells.Item(146, 20).Value = 2
